# The source change being replayed here is purely a cosmetic XML
# attribute-ordering artifact: the upstream repository re-saved this
# fixture after upgrading its OOXML-writing library (Apache POI
# packaging / 3.15 upgrade per the commit message), which made the
# serializer emit element attributes in a different (alphabetical)
# order. Every changed line in the diff carries exactly the same
# attribute names/values as before -- only their left-to-right order
# changed -- and no paragraph text, run, field, style value, or
# document/section property was added, removed, or modified.
#
# Word's object model (and this COM-interop shim) does not expose any
# control over the raw attribute ordering used when an OOXML part is
# serialized -- that is purely an implementation detail of whichever
# library writes the XML, not a document property a macro/automation
# client can observe or set. So there is no Word-level action that
# corresponds to this diff beyond leaving the document's content and
# formatting exactly as they are.
$d = $word.ActiveDocument
